$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 290 (shifts the existing row 290.."354" down to 291.."355")
$ws.Rows.Item(290).Insert()

# Populate the newly inserted row 290 with the new price record
$ws.Range("A290").Value = 3
$ws.Range("B290").Value = "Femacal de La Calera"
$ws.Range("C290").Value = "Coquimbo"
$ws.Range("D290").Value = 45015
$ws.Range("E290").Value = 5
$ws.Range("F290").Value = "Fruta"
$ws.Range("G290").Value = 100101
$ws.Range("H290").Value = "Berries"
$ws.Range("I290").Value = 100101001
$ws.Range("J290").Value = "Arándano (blue)"
$ws.Range("K290").Value = "Sin especificar"
$ws.Range("L290").Value = "Primera"
$ws.Range("M290").Value = 60
$ws.Range("N290").Value = 4000
$ws.Range("O290").Value = 4000
$ws.Range("P290").Value = 4000
$ws.Range("Q290").Value = "$/bandeja 2 kilos"
$ws.Range("R290").Value = "Provincia de Curicó"
$ws.Range("S290").Value = 2000
$ws.Range("T290").Value = 2
